$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Replace the paragraph's existing run text ("Juillet 2022") with the new
# text ("Juillet 2022 " with a trailing space), addressing the run's
# characters directly (not the whole paragraph, which would include the
# trailing paragraph-mark character) so formatting/run structure stays intact.
$oldLen = $para1.Runs(1, 1).Length
$chars = $para1.Characters(1, $oldLen)
$chars.Text = "Juillet 2022 "
